$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("ID Competicao") values for data rows 2-151 were mistakenly
# truncated to 48 and need to be restored to 248.
$ws.Range("B2:B151").Value = 248
